$d = $word.ActiveDocument

# Anchor: the last paragraph in the document currently holds the
# "Student: Sehr gut und schiebt seine Prüfung in den Stapel" text and
# carries the _GoBack bookmark at its very end.
$anchorIndex = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($anchorIndex)
$anchorRange = $anchor.Range

# New paragraphs (SW11 joke block) to insert right after the anchor,
# in document order.
$newTexts = @(
    "",
    "SW11:",
    "Golfer: Arzt, Priester, Ing.",
    "Feuerwehr spielen blind Golf",
    "",
    "Priester: Ich werde sie in mein Gebet nehmen",
    "Arzt: Ich werde einen Kollegen fragen, ob er sie ins Gebet nehmen kann",
    "Ing: Warum können diese nicht nachts spielen?"
)

for ($i = 0; $i -lt $newTexts.Count; $i++) {
    $anchorRange.InsertParagraphAfter()
}

# The freshly inserted paragraphs now occupy indices
# (anchorIndex + 1) .. (anchorIndex + newTexts.Count), in the same
# order they were requested.
for ($i = 0; $i -lt $newTexts.Count; $i++) {
    $p = $d.Paragraphs.Item($anchorIndex + 1 + $i)
    $text = $newTexts[$i]
    if ($text -ne "") {
        $p.Range.Text = $text
    }
}

# Move the _GoBack bookmark from the original anchor paragraph to the
# very end of the new last paragraph, matching the edited document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $lastPara.Range.Duplicate
$endRange.Collapse(0)
$endRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $endRange)

Write-Output "done"
